$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)

$ws.Range("D2").Value = "63.167.28"
$ws.Range("E2").Value = "  -4.36%  "
$ws.Range("D3").Value = "3.280.32"
$ws.Range("E3").Value = "  -6.67%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.605"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.67%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "3.272.79"
$ws.Range("E9").Value = "  -6.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.612"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.156"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000266"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.93%  "
$ws.Range("D15").Value = "3.835.71"
$ws.Range("E15").Value = "  -5.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "17.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.116"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.25%  "
$ws.Range("D18").Value = "3.299.16"
$ws.Range("E18").Value = "  -5.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.37%  "
$ws.Range("D20").Value = "63.089.55"
$ws.Range("E20").Value = "  -4.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.957"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "421.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("E23").Value = "  +4.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.27%  "
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "582.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.105"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.84%  "
$ws.Range("E35").Value = "  -3.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  -6.49%  "
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.86%  "
$ws.Range("D40").Value = "0.0₃0736"
$ws.Range("E40").Value = "  -9.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.359"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").Value = "3.065.18"
$ws.Range("E43").Value = "  -6.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0399"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.58%  "
$ws.Range("E48").Value = "  -3.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.37%  "
